# Adjusting V_reg resistors to use 10k var res
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Remove the old "180 Ohm" resistor row (R1, R2 / RC0603JR-07180RL)
#    -- it is being replaced by the new 910 Ohm part below.
# ---------------------------------------------------------------------
$ws.Rows.Item(6).Delete()

# ---------------------------------------------------------------------
# 2. Insert the new 910 Ohm resistor row (R1, R2) at row 3, copying the
#    formatting of the row below it (the 0 Ohm jumper row) so that the
#    borders / fonts / alignment match the rest of the BOM table.
# ---------------------------------------------------------------------
$ws.Rows.Item(3).Insert()
$ws.Range("B4:H4").Copy()
$ws.Range("B3:H3").PasteSpecial(-4122)
$ws.Range("B3").Value = "CRCW0603910RFKEA"
$ws.Range("C3").Value = "RES, SMD 910 OHM 1% 1/8W 0603"
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = "910 Ohm"
$ws.Range("F3").Value = "R0603"
$ws.Range("G3").Value = "R1, R2"
$ws.Range("H3").Value = "Vishay Dale"

# ---------------------------------------------------------------------
# 3. Insert the new 4.7 kOhm resistor row (R5, R6) at row 5, again
#    copying formatting from the row above (the 0 Ohm jumper row).
# ---------------------------------------------------------------------
$ws.Rows.Item(5).Insert()
$ws.Range("B4:H4").Copy()
$ws.Range("B5:H5").PasteSpecial(-4122)
$ws.Range("B5").Value = "CRCW08054K70FKEAC"
$ws.Range("C5").Value = "RES, 4.7K OHM 1% 1/8W 0805"
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = "4.7 kOhms"
$ws.Range("F5").Value = "R0805"
$ws.Range("G5").Value = "R5, R6"
$ws.Range("H5").Value = "Vishay Dale"

# ---------------------------------------------------------------------
# 4. Tidy up the capacitor value text ("1u" -> "1uF", "2.2u" -> "2.2uF")
# ---------------------------------------------------------------------
$ws.Range("E7").Value = "1uF"
$ws.Range("E8").Value = "2.2uF"

# ---------------------------------------------------------------------
# 5. Widen columns C and H to fit the new, longer text.
# ---------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 35.5
$ws.Columns.Item(8).ColumnWidth = 18.38
